$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts on row 2, row 1 is the header).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$colJ = 10  # column J - "Signal Length [Bit]"
$colK = 11  # column K - "Signal Default"

for ($r = 2; $r -le $lastRow; $r++) {
    $kCell = $ws.Cells.Item($r, $colK)
    $kVal = $kCell.Value2
    if ($kVal -eq $null) {
        # Column K is blank on this row: give it the same look as column J
        # (direct formatting only - named Style is irrelevant here) and fill
        # it with the same placeholder single-space text used by column L.
        $jCell = $ws.Cells.Item($r, $colJ)
        $jCell.Copy()
        $kCell.PasteSpecial(-4122)
        $kCell.Value2 = " "
    }
}
